$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column C ("Förändrad") bumps from 45188 -> 45189 for every existing data row (2..380)
for ($r = 2; $r -le 380; $r++) {
    $ws.Cells.Item($r, 3).Value = 45189
}

# 2) Rows 378 and 379 swap their Beteckning (A), Markägare (F) and Area (G) values.
#    Row 378 was "A 43676-2023" / (no markägare) / 1.1  -> becomes "A 43654-2023" / "Sveaskog" / 0.2
$ws.Cells.Item(378, 1).Value = "A 43654-2023"
$ws.Cells.Item(378, 6).Value = "Sveaskog"
$ws.Cells.Item(378, 7).Value = 0.2

#    Row 379 was "A 43654-2023" / "Sveaskog" / 0.2     -> becomes "A 43676-2023" / (no markägare) / 1.1
$ws.Cells.Item(379, 1).Value = "A 43676-2023"
$ws.Cells.Item(379, 6).ClearContents()
$ws.Cells.Item(379, 7).Value = 1.1

# 3) Row 380 gains an explicit 15pt custom row height (it previously used the sheet default).
$ws.Rows.Item(380).RowHeight = 15

# 4) New row 381 is appended with a fresh record.
$ws.Cells.Item(381, 1).Value = "A 44002-2023"
$ws.Cells.Item(381, 2).Value = 45188
$ws.Cells.Item(381, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(381, 3).Value = 45189
$ws.Cells.Item(381, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(381, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(381, 5).Value = "VAGGERYD"
$ws.Cells.Item(381, 7).Value = 0.9
$ws.Cells.Item(381, 8).Value = 0
$ws.Cells.Item(381, 9).Value = 0
$ws.Cells.Item(381, 10).Value = 0
$ws.Cells.Item(381, 11).Value = 0
$ws.Cells.Item(381, 12).Value = 0
$ws.Cells.Item(381, 13).Value = 0
$ws.Cells.Item(381, 14).Value = 0
$ws.Cells.Item(381, 15).Value = 0
$ws.Cells.Item(381, 16).Value = 0
$ws.Cells.Item(381, 17).Value = 0
$ws.Cells.Item(381, 18).WrapText = $true
